$ppt = $ppt
$p = $ppt.ActivePresentation

$oldDate = "2019/5/1"
$newDate = "2020/7/28"
$oldTypo = "We wan to change the Text Dynamically"
$newTypo = "We want to change the Text Dynamically"

function Update-ShapeText {
    param($shp)

    if (-not $shp.HasTextFrame) { return }
    $tf = $shp.TextFrame
    if (-not $tf.HasText) { return }
    $tr = $tf.TextRange
    $t = $tr.Text
    if ($t -eq $oldDate) {
        $tr.Text = $newDate
    } elseif ($t -eq $oldTypo) {
        $tr.Text = $newTypo
    }
}

function Update-ShapesCollection {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        Update-ShapeText $shapes.Item($i)
    }
}

# 1) Every slide
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Update-ShapesCollection $s.Shapes
}

# 2) Slide master
$master = $p.SlideMaster
Update-ShapesCollection $master.Shapes

# 3) Every slide layout (CustomLayout) under the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-ShapesCollection $layouts.Item($li).Shapes
}

# 4) Notes master - shape TextFrame edits don't persist for the notes
#    master in this host, but the HeadersFooters.DateAndTime object does.
$notesMaster = $p.NotesMaster
$nmDt = $notesMaster.HeadersFooters.DateAndTime
if ($nmDt.Text -eq $oldDate) {
    $nmDt.Text = $newDate
}
Update-ShapesCollection $notesMaster.Shapes
